$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2195.05
$ws.Range("I40").Value = 2961.125
$ws.Range("J40").Value = 1684.3334
$ws.Range("K40").Value = 2961.125
$ws.Range("L40").Value = 1684.3334
$ws.Range("M40").Value = -2786.125
$ws.Range("N40").Value = -2034.3334
$ws.Range("H62").Value = 50665.906
$ws.Range("I62").Value = 86082.75
$ws.Range("J62").Value = 3443.4443
$ws.Range("K62").Value = 86082.75
$ws.Range("L62").Value = 3443.4443
$ws.Range("M62").Value = -85458.75
$ws.Range("N62").Value = -4691.4443
$ws.Range("H64").Value = 428824.34
$ws.Range("I64").Value = 732198.6
$ws.Range("K64").Value = 732198.6
$ws.Range("M64").Value = -731950.6
$ws.Range("H65").Value = 50665.906
$ws.Range("I65").Value = 86082.75
$ws.Range("J65").Value = 3443.4443
$ws.Range("K65").Value = 430413.75
$ws.Range("L65").Value = 17217.2215
$ws.Range("M65").Value = -427293.75
$ws.Range("N65").Value = -23457.2215
$ws.Range("H67").Value = 428824.34
$ws.Range("I67").Value = 732198.6
$ws.Range("K67").Value = 732198.6
$ws.Range("M67").Value = -731340.6
$ws.Range("H69").Value = 83119.766
$ws.Range("I69").Value = 6677.5
$ws.Range("J69").Value = 148641.72
$ws.Range("K69").Value = 20032.5
$ws.Range("L69").Value = 445925.16
$ws.Range("M69").Value = -19158.5
$ws.Range("N69").Value = -447673.16
$ws.Range("H72").Value = 83119.766
$ws.Range("I72").Value = 6677.5
$ws.Range("J72").Value = 148641.72
$ws.Range("K72").Value = 60097.5
$ws.Range("L72").Value = 1337775.48
$ws.Range("M72").Value = -55729.5
$ws.Range("N72").Value = -1346511.48
$ws.Range("H135").Value = 51726190
$ws.Range("I135").Value = 21740896
$ws.Range("J135").Value = 166669820
$ws.Range("K135").Value = 195668064
$ws.Range("L135").Value = 1500028380
$ws.Range("M135").Value = -195665529
$ws.Range("N135").Value = -1500033450
$ws.Range("H137").Value = 2507.034
$ws.Range("I137").Value = 1376.8096
$ws.Range("K137").Value = 4130.4288
$ws.Range("M137").Value = -1580.4288
$ws.Range("H138").Value = 3410.1235
$ws.Range("I138").Value = 1501.3684
$ws.Range("J138").Value = 3928.2144
$ws.Range("K138").Value = 4504.1052
$ws.Range("L138").Value = 11784.6432
$ws.Range("M138").Value = 635.8948
$ws.Range("N138").Value = -22064.6432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7582.904
$ws.Range("I61").Value = 3744.8484
$ws.Range("J61").Value = 14249
$ws.Range("K61").Value = 3744.8484
$ws.Range("L61").Value = 14249
$ws.Range("M61").Value = -3532.8484
$ws.Range("N61").Value = -14673
$ws.Range("H136").Value = 7582.904
$ws.Range("I136").Value = 3744.8484
$ws.Range("J136").Value = 14249
$ws.Range("K136").Value = 11234.5452
$ws.Range("L136").Value = 42747
$ws.Range("M136").Value = -8684.5452
$ws.Range("N136").Value = -47847
$ws.Range("H140").Value = 38042.5
$ws.Range("I140").Value = 29723.334
$ws.Range("J140").Value = 63000
$ws.Range("K140").Value = 29723.334
$ws.Range("L140").Value = 63000
$ws.Range("M140").Value = -24543.334
$ws.Range("N140").Value = -73360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1709.3723
$ws.Range("I31").Value = 1244.629
$ws.Range("J31").Value = 2609.8125
$ws.Range("K31").Value = 1244.629
$ws.Range("L31").Value = 2609.8125
$ws.Range("M31").Value = -949.6289999999999
$ws.Range("N31").Value = -3199.8125
$ws.Range("H34").Value = 1709.3723
$ws.Range("I34").Value = 1244.629
$ws.Range("J34").Value = 2609.8125
$ws.Range("K34").Value = 1244.629
$ws.Range("L34").Value = 2609.8125
$ws.Range("M34").Value = -1042.629
$ws.Range("N34").Value = -3013.8125
$ws.Range("H55").Value = 14285.143
$ws.Range("J55").Value = 14999.333
$ws.Range("L55").Value = 14999.333
$ws.Range("N55").Value = -15629.333
$ws.Range("H62").Value = 2878.889
$ws.Range("I62").Value = 2878.889
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2878.889
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2254.889
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 2878.889
$ws.Range("I65").Value = 2878.889
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 14394.445
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -11274.445
$ws.Range("N65").ClearContents()
$ws.Range("H110").Value = 31140.8
$ws.Range("J110").Value = 31140.8
$ws.Range("L110").Value = 31140.8
$ws.Range("N110").Value = -39320.8
$ws.Range("H112").Value = 70000
$ws.Range("J112").Value = 70000
$ws.Range("L112").Value = 70000
$ws.Range("N112").Value = -72954
$ws.Range("H123").Value = 54220
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 54220
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 54220
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -64020

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1505.2106
$ws.Range("I4").Value = 599.8333
$ws.Range("J4").Value = 1923.0769
$ws.Range("K4").Value = 1799.4999
$ws.Range("L4").Value = 5769.2307
$ws.Range("M4").Value = -1687.4999
$ws.Range("N4").Value = -5993.2307
$ws.Range("H38").Value = 106.53846
$ws.Range("I38").Value = 30
$ws.Range("J38").Value = 129.5
$ws.Range("K38").Value = 90
$ws.Range("L38").Value = 388.5
$ws.Range("M38").Value = 257
$ws.Range("N38").Value = -1082.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 10511
$ws.Range("J18").Value = 10657
$ws.Range("L18").Value = 10657
$ws.Range("N18").Value = -11243
$ws.Range("H70").Value = 5240.6455
$ws.Range("I70").Value = 4851.52
$ws.Range("J70").Value = 5420.7964
$ws.Range("K70").Value = 4851.52
$ws.Range("L70").Value = 5420.7964
$ws.Range("M70").Value = -4581.52
$ws.Range("N70").Value = -5960.7964
$ws.Range("H73").Value = 5240.6455
$ws.Range("I73").Value = 4851.52
$ws.Range("J73").Value = 5420.7964
$ws.Range("K73").Value = 4851.52
$ws.Range("L73").Value = 5420.7964
$ws.Range("M73").Value = -3915.52
$ws.Range("N73").Value = -7292.7964
$ws.Range("H80").Value = 10400.714
$ws.Range("I80").Value = 26402.5
$ws.Range("J80").Value = 4000
$ws.Range("K80").Value = 26402.5
$ws.Range("L80").Value = 4000
$ws.Range("M80").Value = -25404.5
$ws.Range("N80").Value = -5996
$ws.Range("H83").Value = 10400.714
$ws.Range("I83").Value = 26402.5
$ws.Range("J83").Value = 4000
$ws.Range("K83").Value = 132012.5
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = -127020.5
$ws.Range("N83").Value = -29984
$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -34900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3348.2307
$ws.Range("I7").Value = 3034.0588
$ws.Range("J7").Value = 3941.6667
$ws.Range("K7").Value = 3034.0588
$ws.Range("L7").Value = 3941.6667
$ws.Range("M7").Value = -2922.0588
$ws.Range("N7").Value = -4165.6667
$ws.Range("H22").Value = 712.9091
$ws.Range("I22").Value = 911
$ws.Range("J22").Value = 599.7143
$ws.Range("K22").Value = 911
$ws.Range("L22").Value = 599.7143
$ws.Range("M22").Value = -616
$ws.Range("N22").Value = -1189.7143
$ws.Range("H27").Value = 712.9091
$ws.Range("I27").Value = 911
$ws.Range("J27").Value = 599.7143
$ws.Range("K27").Value = 911
$ws.Range("L27").Value = 599.7143
$ws.Range("M27").Value = -804
$ws.Range("N27").Value = -813.7143
$ws.Range("H40").Value = 4570.4
$ws.Range("I40").Value = 4275.5
$ws.Range("J40").Value = 5750
$ws.Range("K40").Value = 4275.5
$ws.Range("L40").Value = 5750
$ws.Range("M40").Value = -4139.5
$ws.Range("N40").Value = -6022
$ws.Range("H55").Value = 250755.06
$ws.Range("I55").Value = 400780.1
$ws.Range("J55").Value = 713.3333
$ws.Range("K55").Value = 400780.1
$ws.Range("L55").Value = 713.3333
$ws.Range("M55").Value = -400607.1
$ws.Range("N55").Value = -1059.3333
$ws.Range("H68").Value = 5950
$ws.Range("I68").Value = 5950
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 5950
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -5201
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 5950
$ws.Range("I71").Value = 5950
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 29750
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -26006
$ws.Range("N71").ClearContents()
$ws.Range("H82").Value = 2399
$ws.Range("I82").Value = 1424
$ws.Range("J82").Value = 3513.2856
$ws.Range("K82").Value = 1424
$ws.Range("L82").Value = 3513.2856
$ws.Range("M82").Value = -1063
$ws.Range("N82").Value = -4235.2856
$ws.Range("H85").Value = 2399
$ws.Range("I85").Value = 1424
$ws.Range("J85").Value = 3513.2856
$ws.Range("K85").Value = 1424
$ws.Range("L85").Value = 3513.2856
$ws.Range("M85").Value = -176
$ws.Range("N85").Value = -6009.2856
$ws.Range("H126").Value = 3348.2307
$ws.Range("I126").Value = 3034.0588
$ws.Range("J126").Value = 3941.6667
$ws.Range("K126").Value = 9102.1764
$ws.Range("L126").Value = 11825.0001
$ws.Range("M126").Value = -6632.1764
$ws.Range("N126").Value = -16765.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3974.2222
$ws.Range("I62").Value = 3875
$ws.Range("J62").Value = 4002.5715
$ws.Range("K62").Value = 3875
$ws.Range("L62").Value = 4002.5715
$ws.Range("M62").Value = -3251
$ws.Range("N62").Value = -5250.5715
$ws.Range("H65").Value = 3974.2222
$ws.Range("I65").Value = 3875
$ws.Range("J65").Value = 4002.5715
$ws.Range("K65").Value = 19375
$ws.Range("L65").Value = 20012.8575
$ws.Range("M65").Value = -16255
$ws.Range("N65").Value = -26252.8575
$ws.Range("H68").Value = 36271
$ws.Range("J68").Value = 36271
$ws.Range("L68").Value = 36271
$ws.Range("N68").Value = -37893
$ws.Range("H71").Value = 36271
$ws.Range("J71").Value = 36271
$ws.Range("L71").Value = 108813
$ws.Range("N71").Value = -116925
$ws.Range("H121").Value = 36113.332
$ws.Range("J121").Value = 36113.332
$ws.Range("L121").Value = 36113.332
$ws.Range("N121").Value = -39607.332
$ws.Range("H132").Value = 2575.9788
$ws.Range("I132").Value = 2058.5454
$ws.Range("J132").Value = 3795.6428
$ws.Range("K132").Value = 6175.6362
$ws.Range("L132").Value = 11386.9284
$ws.Range("M132").Value = -3645.6362
$ws.Range("N132").Value = -16446.9284
